$p = $ppt.ActivePresentation

# Slide 1: Title "First" + " " + "slide" -> "First " + "slide"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$sub1 = $tr1.Characters(1, 6)
$sub1.Text = "First "

# Slide 3: Title "Third" + " " + "slide" -> "Third " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$sub3 = $tr3.Characters(1, 6)
$sub3.Text = "Third "
